$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, pushing existing rows 145..224 down to 146..225.
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new record's data.
$ws.Cells.Item(145, 1).Value = 3
$ws.Cells.Item(145, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(145, 3).Value = "Coquimbo"
$ws.Cells.Item(145, 4).Value = 44488
$ws.Cells.Item(145, 5).Value = 5
$ws.Cells.Item(145, 6).Value = 100112040
$ws.Cells.Item(145, 7).Value = "Cilantro"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 160
$ws.Cells.Item(145, 11).Value = 2500
$ws.Cells.Item(145, 12).Value = 2500
$ws.Cells.Item(145, 13).Value = 2500
$ws.Cells.Item(145, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(145, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(145, 16).Value = 833
$ws.Cells.Item(145, 17).Value = 3
$ws.Cells.Item(145, 18).Value = "Hortaliza"
